$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.698109865188599
$ws.Range("B1").Value = 6.696359157562256
$ws.Range("C1").Value = 2.797043323516846
$ws.Range("D1").Value = 1.569325089454651
$ws.Range("E1").Value = 1.21110463142395
